$p = $ppt.ActivePresentation

# The deck ships two DrawingML theme parts: theme1.xml ("Office Theme",
# only wired up via the Notes Master) and theme2.xml ("Integral", the
# theme actually driving the Slide Master / every slide layout). The
# edit swaps the two themes' color schemes (their font/format schemes
# are already identical) so the deck's visible design becomes the
# previously-dormant "Office Theme" palette.
#
# The Slide Master's ColorScheme is the live, rendered theme palette -
# write the target ("Office Theme") colors into it, in clrScheme slot
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$slideMaster = $p.SlideMaster
$cs = $slideMaster.ColorScheme

$cs.Colors(1).RGB  = 0x000000   # dk1
$cs.Colors(2).RGB  = 0xFFFFFF   # lt1
$cs.Colors(3).RGB  = 0x6A5444   # dk2      (44546A)
$cs.Colors(4).RGB  = 0xE6E6E7   # lt2      (E7E6E6)
$cs.Colors(5).RGB  = 0xD59B5B   # accent1  (5B9BD5)
$cs.Colors(6).RGB  = 0x317DED   # accent2  (ED7D31)
$cs.Colors(7).RGB  = 0xA5A5A5   # accent3  (A5A5A5)
$cs.Colors(8).RGB  = 0x00C0FF   # accent4  (FFC000)
$cs.Colors(9).RGB  = 0xC47244   # accent5  (4472C4)
$cs.Colors(10).RGB = 0x47AD70   # accent6  (70AD47)
$cs.Colors(11).RGB = 0xC16305   # hlink    (0563C1)
$cs.Colors(12).RGB = 0x724F95   # folHlink (954F72)
